# Update - Suppression de la colonne "Disponible à la vente" du fichier d'import
#
# The "Catalogue" sheet has a table (Tableau4) whose column J is
# "Disponible à la vente (si "Non", les consommateurs ne peuvent pas
# commander le produit)". That whole column is removed: the cells shift
# left (the old "Description" column K becomes J) and the table shrinks
# from 11 to 10 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catalogue")
$lo = $ws.ListObjects.Item(1)

# Prime the header text that should end up in the last remaining table
# column (the "Description" column sliding from K into J) so the table's
# column-name bookkeeping picks it up when we resize below.
$ws.Range("J1").Value = "Description"

# Remove the whole column; everything to the right shifts one column left.
$ws.Range("J1").EntireColumn.Delete()

# Shrink the table definition to match the new, narrower extent.
$lo.Resize($ws.Range("A1:J1048576"))

# The header row auto-fit to a shorter height now that one of the two
# wrapped, two-line headers is gone.
$ws.Rows.Item(1).RowHeight = 60.75

# Match the post-edit selection (whole column J selected).
$ws.Activate() | Out-Null
$ws.Range("J1:J1048576").Select() | Out-Null
